# Trade #8 closed at 2026-02-17 04:06:41 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet updates ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.05      # Total P&L %
$summary.Range("B6").Value = 8         # Total Trades
$summary.Range("B7").Value = 3         # Winning Trades
$summary.Range("B9").Value = 37.5      # Win Rate %

# ---- Strategy Status sheet updates (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 8          # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 37.5       # Win Rate %

# ---- New trade row (row 9) appended to "All Trades" and "MarketMaking" sheets ----
$newRow = @(8, "2026-02-17", "04:06:35", "MarketMaking", "DOWN", 0.78, 0.83, "CLOSED", 6.4103, 0.05, 100.02, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.12)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $cell = $ws.Cells.Item(9, $col)
        if ($col -eq 2) {
            # "2026-02-17" looks like a date; force it to stay plain text
            # (matching the rest of the sheet) instead of being auto-converted
            # into a date serial number by Excel's smart input parsing.
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$col - 1]
            $cell.ClearFormats()
        } else {
            $cell.Value = $newRow[$col - 1]
        }
    }
}
